$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cyclic content shift across rows 55-59:
#   before-row 59 -> row 55
#   before-row 55 -> row 57
#   before-row 57 -> row 59
#   before-row 58 -> row 56
#   before-row 56 -> row 58


# Row 55
$ws.Range("A55").Value = 111898889
$ws.Range("B55").Value = 98535
$ws.Range("D55").Value = "LC"
$ws.Range("E55").Value = 222498
$ws.Range("F55").Value = "Blåsippa"
$ws.Range("G55").Value = "Hepatica nobilis"
$ws.Range("H55").Value = "Schreb."
$ws.Range("K55").Value = "fullt utvecklade blad"
$ws.Range("L55").Value = ""
$ws.Range("Q55").Value = 650135.0421630922
$ws.Range("R55").Value = 6654002.501842719
$ws.Range("AH55").Value = "Ängsbarrskog"
$ws.Range("AI55").Value = "Ungskog"
$ws.Range("AJ55").Value = ""
$ws.Range("AK55").Value = ""
$ws.Range("AM55").Value = ""
$ws.Range("AO55").Value = ""

# Row 56
$ws.Range("A56").Value = 111898660
$ws.Range("B56").Value = 100532
$ws.Range("D56").Value = "CR"
$ws.Range("E56").Value = 223246
$ws.Range("F56").Value = "Skogsalm"
$ws.Range("G56").Value = "Ulmus glabra"
$ws.Range("H56").Value = "Huds."
$ws.Range("L56").Value = ""
$ws.Range("Q56").Value = 650054.1336129439
$ws.Range("R56").Value = 6654018.240072312
$ws.Range("AC56").Value = "Stammens omkrets i brösthöjd: 64 cm"
$ws.Range("AJ56").Value = ""
$ws.Range("AK56").Value = ""
$ws.Range("AM56").Value = ""
$ws.Range("AO56").Value = ""

# Row 57
$ws.Range("A57").Value = 111898336
$ws.Range("B57").Value = 89405
$ws.Range("D57").Value = "NT"
$ws.Range("E57").Value = 1202
$ws.Range("F57").Value = "Ullticka"
$ws.Range("G57").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H57").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I57").Value = ""
$ws.Range("J57").Value = ""
$ws.Range("Q57").Value = 650105.085176448
$ws.Range("R57").Value = 6654011.298884101
$ws.Range("AH57").Value = "Ängsblandskog"
$ws.Range("AI57").Value = ""
$ws.Range("AJ57").Value = "gran"
$ws.Range("AK57").Value = "Picea abies"
$ws.Range("AM57").Value = "Liggande död trädstam, utan markontakt"
$ws.Range("AO57").Value = "Horizontal, dead without ground contact # Picea abies"

# Row 58
$ws.Range("A58").Value = 111898507
$ws.Range("B58").Value = 89845
$ws.Range("D58").Value = "VU"
$ws.Range("E58").Value = 1209
$ws.Range("F58").Value = "Rynkskinn"
$ws.Range("G58").Value = "Phlebia centrifuga"
$ws.Range("H58").Value = "P.Karst."
$ws.Range("L58").Value = ""
$ws.Range("Q58").Value = 650086.8716060545
$ws.Range("R58").Value = 6654015.064976334
$ws.Range("AC58").Value = ""
$ws.Range("AJ58").Value = "gran"
$ws.Range("AK58").Value = "Picea abies"
$ws.Range("AM58").Value = "Liggande död trädstam, utan markontakt"
$ws.Range("AO58").Value = "Horizontal, dead without ground contact # Picea abies"

# Row 59
$ws.Range("A59").Value = 111898191
$ws.Range("B59").Value = 90332
$ws.Range("E59").Value = 4769
$ws.Range("F59").Value = "Svavelriska"
$ws.Range("G59").Value = "Lactarius scrobiculatus"
$ws.Range("H59").Value = "(Scop.:Fr.) Fr."
$ws.Range("I59").Value = "'2"
$ws.Range("J59").Value = "fruktkroppar"
$ws.Range("K59").Value = ""
$ws.Range("L59").Value = ""
